$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the ": NoSettings" tag suffix to ": test" across all the
#    "Process Emissions before CCS[...]" / "Industrial Sector Energy Related
#    Emissions before CCS[...]" row labels on the "BAU Emissions" sheet.
#    These labels live in the shared-string table, so a single UsedRange
#    Replace takes care of every occurrence at once.
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
[void]$wsBau.UsedRange.Replace(" : NoSettings", " : test", -4163, 1, $false, $false, $true)

# ---------------------------------------------------------------------------
# 2) Update the updated model-run numbers on row 94 (Industrial Sector
#    Energy Related Emissions before CCS[natural gas if,iron and steel
#    241,CO2]) for years 2032-2050 (columns M:AE).
# ---------------------------------------------------------------------------
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94").Value = 5005380
$wsBau.Range("R94").Value = 5005380
$wsBau.Range("S94").Value = 5005380
$wsBau.Range("T94").Value = 5005380
$wsBau.Range("U94").Value = 5005380
$wsBau.Range("V94").Value = 5005380
$wsBau.Range("W94").Value = 5005380
$wsBau.Range("X94").Value = 5005380
$wsBau.Range("Y94").Value = 5005380
$wsBau.Range("Z94").Value = 5005380
$wsBau.Range("AA94").Value = 5005380
$wsBau.Range("AB94").Value = 5005380
$wsBau.Range("AC94").Value = 5005380
$wsBau.Range("AD94").Value = 5005380
$wsBau.Range("AE94").Value = 5005380

# ---------------------------------------------------------------------------
# 3) Bump the "last updated" serial date stamp on the About sheet.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# ---------------------------------------------------------------------------
# 4) Update the view state: scroll/select the BAU Emissions sheet to the
#    new region, then finish with the About sheet active/selected (matches
#    the saved workbook view).
# ---------------------------------------------------------------------------
$wsBau.Activate()
[void]$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()
[void]$wsAbout.Range("E29").Select()
